$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "#ParoNacional4J"
$ws.Range("B6").Value = 32000
